# Apply updated loading-percent values ("case with 380 kV done")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell = "B2"; Value = 9.077415285415968},
    @{Cell = "C2"; Value = 5.417613208607194},
    @{Cell = "D2"; Value = 4.992241250410425},
    @{Cell = "F2"; Value = 22.04533146417304},
    @{Cell = "G2"; Value = 3.625120860072542},
    @{Cell = "I2"; Value = 19.16784200727609},
    @{Cell = "K2"; Value = 8.683681599227008},
    @{Cell = "O2"; Value = 19.98374392759137},
    @{Cell = "B3"; Value = 8.688384084678983},
    @{Cell = "C3"; Value = 5.211784515551813},
    @{Cell = "D3"; Value = 4.890931008010105},
    @{Cell = "F3"; Value = 22.16078855855167},
    @{Cell = "G3"; Value = 3.626581312923057},
    @{Cell = "I3"; Value = 19.29197477008969},
    @{Cell = "K3"; Value = 8.405505909228753},
    @{Cell = "O3"; Value = 20.10546929477558},
    @{Cell = "B4"; Value = 8.440656668651931},
    @{Cell = "C4"; Value = 5.080398039259375},
    @{Cell = "D4"; Value = 4.827026935030815},
    @{Cell = "F4"; Value = 22.2382985536303},
    @{Cell = "G4"; Value = 3.627525325314044},
    @{Cell = "I4"; Value = 19.37226854679644},
    @{Cell = "K4"; Value = 8.230261231717664},
    @{Cell = "O4"; Value = 20.18491667291989},
    @{Cell = "B5"; Value = 8.337608438107596},
    @{Cell = "C5"; Value = 5.025652410229435},
    @{Cell = "D5"; Value = 4.800579556237909},
    @{Cell = "F5"; Value = 22.27154337451682},
    @{Cell = "G5"; Value = 3.627921944294474},
    @{Cell = "I5"; Value = 19.40601579064269},
    @{Cell = "K5"; Value = 8.157828602906807},
    @{Cell = "O5"; Value = 20.21847571274299},
    @{Cell = "B6"; Value = 8.320374858171556},
    @{Cell = "C6"; Value = 5.016490910181251},
    @{Cell = "D6"; Value = 4.796164115047523},
    @{Cell = "F6"; Value = 22.27716368128531},
    @{Cell = "G6"; Value = 3.62798852392797},
    @{Cell = "I6"; Value = 19.41168157484524},
    @{Cell = "K6"; Value = 8.145742715803264},
    @{Cell = "O6"; Value = 20.22411964837359},
    @{Cell = "B7"; Value = 8.439275230206809},
    @{Cell = "C7"; Value = 5.07966452097207},
    @{Cell = "D7"; Value = 4.826671871150872},
    @{Cell = "F7"; Value = 22.23874019552542},
    @{Cell = "G7"; Value = 3.62753062592051},
    @{Cell = "I7"; Value = 19.3727195134693},
    @{Cell = "K7"; Value = 8.229288369641656},
    @{Cell = "O7"; Value = 20.18536446971673},
    @{Cell = "B8"; Value = 8.94518946394005},
    @{Cell = "C8"; Value = 5.347713253245643},
    @{Cell = "D8"; Value = 4.9576745910147},
    @{Cell = "F8"; Value = 22.08376403849484},
    @{Cell = "G8"; Value = 3.625614630733565},
    @{Cell = "I8"; Value = 19.20979826141429},
    @{Cell = "K8"; Value = 8.588745334039091},
    @{Cell = "O8"; Value = 20.02473768829067},
    @{Cell = "B9"; Value = 9.86208630346759},
    @{Cell = "C9"; Value = 5.831564955453903},
    @{Cell = "D9"; Value = 5.200209545771432},
    @{Cell = "F9"; Value = 21.83261673269427},
    @{Cell = "G9"; Value = 3.62223097059041},
    @{Cell = "I9"; Value = 18.92255780172722},
    @{Cell = "K9"; Value = 9.254484889629545},
    @{Cell = "O9"; Value = 19.74710817889208},
    @{Cell = "B10"; Value = 10.48439808464353},
    @{Cell = "C10"; Value = 6.159225772462222},
    @{Cell = "D10"; Value = 5.368551873140667},
    @{Cell = "F10"; Value = 21.68058981534267},
    @{Cell = "G10"; Value = 3.619970488565501},
    @{Cell = "I10"; Value = 18.73105211718832},
    @{Cell = "K10"; Value = 9.715101232681896},
    @{Cell = "O10"; Value = 19.5659193180783},
    @{Cell = "B11"; Value = 10.75546385034703},
    @{Cell = "C11"; Value = 6.301854386377939},
    @{Cell = "D11"; Value = 5.442789125697383},
    @{Cell = "F11"; Value = 21.61855222658457},
    @{Cell = "G11"; Value = 3.618990618860016},
    @{Cell = "I11"; Value = 18.64814502344577},
    @{Cell = "K11"; Value = 9.917638197675698},
    @{Cell = "O11"; Value = 19.4884432983011},
    @{Cell = "B12"; Value = 10.85632024674157},
    @{Cell = "C12"; Value = 6.354913436178626},
    @{Cell = "D12"; Value = 5.470548195313449},
    @{Cell = "F12"; Value = 21.59608960935995},
    @{Cell = "G12"; Value = 3.618626496346418},
    @{Cell = "I12"; Value = 18.61735373389052},
    @{Cell = "K12"; Value = 9.993268626661859},
    @{Cell = "O12"; Value = 19.4598172570648},
    @{Cell = "B13"; Value = 10.83467951013199},
    @{Cell = "C13"; Value = 6.343528922497646},
    @{Cell = "D13"; Value = 5.46458572528541},
    @{Cell = "F13"; Value = 21.60088145015417},
    @{Cell = "G13"; Value = 3.618704608810027},
    @{Cell = "I13"; Value = 18.62395835757791},
    @{Cell = "K13"; Value = 9.977028540994715},
    @{Cell = "O13"; Value = 19.46595068459122},
    @{Cell = "B14"; Value = 10.76379756214503},
    @{Cell = "C14"; Value = 6.306238794141769},
    @{Cell = "D14"; Value = 5.445080043134146},
    @{Cell = "F14"; Value = 21.61668354468846},
    @{Cell = "G14"; Value = 3.618960523518695},
    @{Cell = "I14"; Value = 18.64559971254651},
    @{Cell = "K14"; Value = 9.923882037173783},
    @{Cell = "O14"; Value = 19.48607393230074},
    @{Cell = "B15"; Value = 10.72014552703362},
    @{Cell = "C15"; Value = 6.283272883629512},
    @{Cell = "D15"; Value = 5.4330858248761},
    @{Cell = "F15"; Value = 21.62649704714347},
    @{Cell = "G15"; Value = 3.619118180587567},
    @{Cell = "I15"; Value = 18.65893425874711},
    @{Cell = "K15"; Value = 9.891187808671594},
    @{Cell = "O15"; Value = 19.49849281441977},
    @{Cell = "B16"; Value = 10.4664358598945},
    @{Cell = "C16"; Value = 6.149772927717372},
    @{Cell = "D16"; Value = 5.363651751896057},
    @{Cell = "F16"; Value = 21.68478789276237},
    @{Cell = "G16"; Value = 3.620035497247455},
    @{Cell = "I16"; Value = 18.73655487022118},
    @{Cell = "K16"; Value = 9.701718880772773},
    @{Cell = "O16"; Value = 19.57108219973793},
    @{Cell = "B17"; Value = 10.3076663484499},
    @{Cell = "C17"; Value = 6.066208494785764},
    @{Cell = "D17"; Value = 5.320444377783287},
    @{Cell = "F17"; Value = 21.72237525897404},
    @{Cell = "G17"; Value = 3.620610623789851},
    @{Cell = "I17"; Value = 18.78524972344603},
    @{Cell = "K17"; Value = 9.583648586446628},
    @{Cell = "O17"; Value = 19.61688132607458},
    @{Cell = "B18"; Value = 10.21521725806912},
    @{Cell = "C18"; Value = 6.017540960797225},
    @{Cell = "D18"; Value = 5.295373037544786},
    @{Cell = "F18"; Value = 21.74466436781039},
    @{Cell = "G18"; Value = 3.620945982183458},
    @{Cell = "I18"; Value = 18.81365405044004},
    @{Cell = "K18"; Value = 9.515081922927667},
    @{Cell = "O18"; Value = 19.64368940973886},
    @{Cell = "B19"; Value = 10.18372368195302},
    @{Cell = "C19"; Value = 6.00096021869991},
    @{Cell = "D19"; Value = 5.286847106322748},
    @{Cell = "F19"; Value = 21.75232594587691},
    @{Cell = "G19"; Value = 3.621060313078064},
    @{Cell = "I19"; Value = 18.82333938613533},
    @{Cell = "K19"; Value = 9.491755724993411},
    @{Cell = "O19"; Value = 19.65284611173485},
    @{Cell = "B20"; Value = 10.32468495978658},
    @{Cell = "C20"; Value = 6.075166739727234},
    @{Cell = "D20"; Value = 5.325066724943891},
    @{Cell = "F20"; Value = 21.71830465035197},
    @{Cell = "G20"; Value = 3.620548928792378},
    @{Cell = "I20"; Value = 18.78002506525661},
    @{Cell = "K20"; Value = 9.596285690670895},
    @{Cell = "O20"; Value = 19.61195773238001},
    @{Cell = "B21"; Value = 10.78466633495506},
    @{Cell = "C21"; Value = 6.317217824751483},
    @{Cell = "D21"; Value = 5.450819036775491},
    @{Cell = "F21"; Value = 21.61201409598118},
    @{Cell = "G21"; Value = 3.618885167263264},
    @{Cell = "I21"; Value = 18.63922674629473},
    @{Cell = "K21"; Value = 9.939521821925307},
    @{Cell = "O21"; Value = 19.48014390289089},
    @{Cell = "B22"; Value = 11.07483364361409},
    @{Cell = "C22"; Value = 6.469857289878593},
    @{Cell = "D22"; Value = 5.530941437498027},
    @{Cell = "F22"; Value = 21.54855137805193},
    @{Cell = "G22"; Value = 3.617838198388311},
    @{Cell = "I22"; Value = 18.55072579683506},
    @{Cell = "K22"; Value = 10.17663676390809},
    @{Cell = "O22"; Value = 19.39814913937109},
    @{Cell = "B23"; Value = 10.92094040136168},
    @{Cell = "C23"; Value = 6.388907014231108},
    @{Cell = "D23"; Value = 5.488372479739074},
    @{Cell = "F23"; Value = 21.58187135161836},
    @{Cell = "G23"; Value = 3.618393299695633},
    @{Cell = "I23"; Value = 18.59763893523012},
    @{Cell = "K23"; Value = 10.04435695556056},
    @{Cell = "O23"; Value = 19.44153092629877},
    @{Cell = "B24"; Value = 10.31699448637187},
    @{Cell = "C24"; Value = 6.071118662546496},
    @{Cell = "D24"; Value = 5.322977679592128},
    @{Cell = "F24"; Value = 21.72014285596302},
    @{Cell = "G24"; Value = 3.620576806425432},
    @{Cell = "I24"; Value = 18.78238585907598},
    @{Cell = "K24"; Value = 9.59057458950134},
    @{Cell = "O24"; Value = 19.6141822016076},
    @{Cell = "B25"; Value = 9.622712768124526},
    @{Cell = "C25"; Value = 5.70540794094848},
    @{Cell = "D25"; Value = 5.136252337007741},
    @{Cell = "F25"; Value = 21.89487627632886},
    @{Cell = "G25"; Value = 3.623106578827441},
    @{Cell = "I25"; Value = 18.99682514550969},
    @{Cell = "K25"; Value = 9.079086246075715},
    @{Cell = "O25"; Value = 19.81821399725273}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
